$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82-88 down to 83-89
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new record's data
$ws.Cells.Item(82, 1).Value = 5
$ws.Cells.Item(82, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value = "Maule"
$ws.Cells.Item(82, 4).Value = 44461
$ws.Cells.Item(82, 5).Value = 7
$ws.Cells.Item(82, 6).Value = 100112031
$ws.Cells.Item(82, 7).Value = "Poroto verde"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 100
$ws.Cells.Item(82, 11).Value = 32000
$ws.Cells.Item(82, 12).Value = 32000
$ws.Cells.Item(82, 13).Value = 32000
$ws.Cells.Item(82, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 1280
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
